$d = $word.ActiveDocument

$d.Content.Find.Execute('2025-02-25 Tuesday', $true, $false, $false, $false, $false, $true, 1, $false, '2025-02-26 Wednesday', 2) | Out-Null
$d.Content.Find.Execute('95×98=9310', $true, $false, $false, $false, $false, $true, 1, $false, '28×75=2100', 2) | Out-Null
$d.Content.Find.Execute('89×83=7387', $true, $false, $false, $false, $false, $true, 1, $false, '80×39=3120', 2) | Out-Null
$d.Content.Find.Execute('50×20=1000', $true, $false, $false, $false, $false, $true, 1, $false, '84×70=5880', 2) | Out-Null
$d.Content.Find.Execute('97×93=9021', $true, $false, $false, $false, $false, $true, 1, $false, '50×78=3900', 2) | Out-Null
$d.Content.Find.Execute('55×60=3300', $true, $false, $false, $false, $false, $true, 1, $false, '98×32=3136', 2) | Out-Null
$d.Content.Find.Execute('96×99=9504', $true, $false, $false, $false, $false, $true, 1, $false, '39×91=3549', 2) | Out-Null
$d.Content.Find.Execute('36×97=3492', $true, $false, $false, $false, $false, $true, 1, $false, '27×32=864', 2) | Out-Null
$d.Content.Find.Execute('68×68=4624', $true, $false, $false, $false, $false, $true, 1, $false, '74×89=6586', 2) | Out-Null
$d.Content.Find.Execute('72×43=3096', $true, $false, $false, $false, $false, $true, 1, $false, '63×14=882', 2) | Out-Null
$d.Content.Find.Execute('65×70=4550', $true, $false, $false, $false, $false, $true, 1, $false, '20×64=1280', 2) | Out-Null
$d.Content.Find.Execute('96×43=4128', $true, $false, $false, $false, $false, $true, 1, $false, '69×54=3726', 2) | Out-Null
$d.Content.Find.Execute('45×64=2880', $true, $false, $false, $false, $false, $true, 1, $false, '82×20=1640', 2) | Out-Null
$d.Content.Find.Execute('50×37=1850', $true, $false, $false, $false, $false, $true, 1, $false, '98×94=9212', 2) | Out-Null
$d.Content.Find.Execute('95×33=3135', $true, $false, $false, $false, $false, $true, 1, $false, '57×82=4674', 2) | Out-Null
$d.Content.Find.Execute('39×71=2769', $true, $false, $false, $false, $false, $true, 1, $false, '36×46=1656', 2) | Out-Null
$d.Content.Find.Execute('49×95=4655', $true, $false, $false, $false, $false, $true, 1, $false, '35×64=2240', 2) | Out-Null
$d.Content.Find.Execute('86×40=3440', $true, $false, $false, $false, $false, $true, 1, $false, '78×23=1794', 2) | Out-Null
$d.Content.Find.Execute('38×26=988', $true, $false, $false, $false, $false, $true, 1, $false, '49×92=4508', 2) | Out-Null
$d.Content.Find.Execute('41×78=3198', $true, $false, $false, $false, $false, $true, 1, $false, '64×26=1664', 2) | Out-Null
$d.Content.Find.Execute('99×57=5643', $true, $false, $false, $false, $false, $true, 1, $false, '29×99=2871', 2) | Out-Null
$d.Content.Find.Execute('45×91=4095', $true, $false, $false, $false, $false, $true, 1, $false, '98×36=3528', 2) | Out-Null
$d.Content.Find.Execute('18×31=558', $true, $false, $false, $false, $false, $true, 1, $false, '16×83=1328', 2) | Out-Null
$d.Content.Find.Execute('14×61=854', $true, $false, $false, $false, $false, $true, 1, $false, '11×53=583', 2) | Out-Null
$d.Content.Find.Execute('93×69=6417', $true, $false, $false, $false, $false, $true, 1, $false, '47×92=4324', 2) | Out-Null
$d.Content.Find.Execute('54×89=4806', $true, $false, $false, $false, $false, $true, 1, $false, '95×55=5225', 2) | Out-Null

$d.Save()
